# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.466.69'

$ws.Range("D3").Value = '2.092.18'
$ws.Range("E3").Value = '  -1.32%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("E7").Value = '  -0.75%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4406'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.78'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +14.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08922'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.66%  '

$ws.Range("E11").Value = '  -3.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.34%  '

$ws.Range("D13").Value = '2.095.22'
$ws.Range("E13").Value = '  -1.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.685'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.680'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.26%  '

$ws.Range("E16").Value = '  -2.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001122'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06597'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.59%  '

$ws.Range("E20").Value = '  -0.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.262'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.33%  '

$ws.Range("D23").Value = '30.498.83'
$ws.Range("E23").Value = '  -1.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.31'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.322'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.54%  '

$ws.Range("D26").Value = '2.336.99'
$ws.Range("E26").Value = '  -1.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.21%  '

$ws.Range("E28").Value = '  -1.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.186'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1069'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.661'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.94%  '

$ws.Range("E34").Value = '  -1.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.898'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.00%  '

$ws.Range("E36").Value = '  +4.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02566'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06820'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.03%  '

$ws.Range("E39").Value = '  -2.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6881'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.253'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6340'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.197'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.627'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.231'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.29%  '

$ws.Range("E50").Value = '  -3.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.70%  '
